$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cell values updated per the crypto-price refresh diff.
# D-column price cells must stay plain text (values like "1.001" or
# "303.68" would otherwise be auto-detected as numbers by Excel), so we
# force the Text number format before assigning, then restore the default
# "Normal" style so no stray formatting is left behind.

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '23.241.29'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  +1.02%  '

$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.602.99'
$ws.Range('D3').Style = 'Normal'

$ws.Range('E4').Value = '  -0.14%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '1.001'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -0.10%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '303.68'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +0.89%  '

$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.3768'
$ws.Range('D7').Style = 'Normal'

$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '51.84'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +5.15%  '

$ws.Range('E9').Value = '  +0.23%  '

$ws.Range('E10').Value = '  +1.40%  '

$ws.Range('E11').Value = '  -0.13%  '

$ws.Range('E12').Value = '  +0.27%  '

$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '22.81'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +0.13%  '

$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '6.596'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +0.30%  '

$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '7.425'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +0.33%  '

$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.00001250'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +0.64%  '

$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '1.607.89'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +0.63%  '

$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '94.05'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +2.22%  '

$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.06922'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +0.71%  '

$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '18.18'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -0.23%  '

$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '6.531'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -0.25%  '

$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '1.001'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -0.21%  '

$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '12.95'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -1.53%  '

$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '23.246.66'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +1.01%  '

$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '3.036'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +8.80%  '

$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '2.381'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +0.72%  '

$ws.Range('E27').Value = '  +0.83%  '

$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '150.09'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -0.22%  '

$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '5.255'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +0.02%  '

$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '134.66'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +0.71%  '

$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '2.401'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +4.09%  '

$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '6.754'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -1.13%  '

$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '1.782.51'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -0.09%  '

$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.9615'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +0.06%  '

$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.07503'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -1.59%  '

$ws.Range('B36').Value = 'VeChain'
$ws.Range('C36').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.02748'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +1.48%  '

$ws.Range('B37').Value = 'FraxShare'
$ws.Range('C37').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '10.32'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -0.20%  '

$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.2537'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +0.22%  '

$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '6.126'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -2.31%  '

$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.08828'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -0.25%  '

$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '1.398'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +2.39%  '

$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.7114'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +0.94%  '

$ws.Range('E43').Value = '  -0.03%  '

$ws.Range('E44').Value = '  +2.04%  '

$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.6549'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -0.88%  '

$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '2.318'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +0.29%  '

$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.9998'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -0.11%  '

$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '4.017'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +0.69%  '

$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '132.70'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +0.26%  '

$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.07951'
$ws.Range('D50').Style = 'Normal'

$ws.Range('E51').Value = '  -1.87%  '
